$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = "rana.abozaid@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G4").Value = "aya.saeed@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G5").Value = "aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"
$ws.Range("G14").Value = "nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G15").Value = "nourhan.mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg"
$ws.Range("G17").Value = "youstina.magdy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G18").Value = "neveen.nashaat@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G19").Value = "rana.abozaid@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G21").Value = "aya.saeed@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G22").Value = "aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"
$ws.Range("G31").Value = "nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G32").Value = "nourhan.mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg"
$ws.Range("G34").Value = "youstina.magdy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G35").Value = "neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"
$ws.Range("G37").Value = "Kerelos.zareef@med.asu.edu.eg, System"
$ws.Range("G45").Value = "Rania.a.youssef@med.asu.edu.eg, backup@backdoor.com, mohamed.saleem@med.asu.edu.eg, System"
$ws.Range("G48").Value = "afaf.abdallah@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G49").Value = "nourhan.mostafa@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg"
$ws.Range("G51").Value = "neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("G52").Value = "naema.gomaa@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Range("G54").Value = "Kerelos.zareef@med.asu.edu.eg, System"
$ws.Range("G62").Value = "Rania.a.youssef@med.asu.edu.eg, backup@backdoor.com, mohamed.saleem@med.asu.edu.eg, System"
$ws.Range("G65").Value = "afaf.abdallah@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G66").Value = "nourhan.mostafa@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg"
$ws.Range("G68").Value = "neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("G69").Value = "naema.gomaa@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Range("G70").Value = "rana.abozaid@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G72").Value = "wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G73").Value = "aya.saeed@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G76").Value = "Rania.a.youssef@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G83").Value = "afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg, Youstina.ibrahim@med.asu.edu.eg"
$ws.Range("G85").Value = "neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G86").Value = "neveen.nashaat@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G87").Value = "rana.abozaid@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G89").Value = "aya.saeed@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G98").Value = "nourhanmohamed@med.asu.edu.eg, user@user.com, nourhanhosni@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G100").Value = "nourhanmohamed@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg"
$ws.Range("G102").Value = "neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G103").Value = "neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Range("G105").Value = "backup@backdoor.com, System"
$ws.Range("G106").Value = "aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"
$ws.Range("G116").Value = "afaf.abdallah@med.asu.edu.eg, enas.omran@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("G117").Value = "Amr-Saeed@med.asu.edu.eg`n, enas.omran@med.asu.edu.eg, user@user.com, takyas@med.asu.edu.eg"
$ws.Range("G119").Value = "neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, marinasorial@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("G120").Value = "neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Range("G122").Value = "backup@backdoor.com, System"
$ws.Range("G123").Value = "aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"
$ws.Range("G133").Value = "afaf.abdallah@med.asu.edu.eg, enas.omran@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("G134").Value = "Amr-Saeed@med.asu.edu.eg`n, enas.omran@med.asu.edu.eg, user@user.com, takyas@med.asu.edu.eg"
$ws.Range("G136").Value = "neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, marinasorial@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("G137").Value = "neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Range("G144").Value = "nourhanmamdouh@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G150").Value = "afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg, Youstina.ibrahim@med.asu.edu.eg"
$ws.Range("G151").Value = "nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg, System, randa.rabea@med.asu.edu.eg"
